$wb = $excel.ActiveWorkbook

# --- Settings sheet: add the System1 login URL as a hyperlink in B6 ---
$settings = $wb.Worksheets.Item("Settings")
$settings.Hyperlinks.Add($settings.Range("B6"), "https://acme-test.uipath.com/login") | Out-Null

# --- Constants sheet: add Timeout* constants (rows 12-16) ---
$constants = $wb.Worksheets.Item("Constants")
$constants.Range("A12").Value = "TimeoutXS"
$constants.Range("B12").Value = 1000
$constants.Range("A13").Value = "TimeoutS"
$constants.Range("B13").Value = 5000
$constants.Range("A14").Value = "TimeoutM"
$constants.Range("B14").Value = 10000
$constants.Range("A15").Value = "TimeoutL"
$constants.Range("B15").Value = 30000
$constants.Range("A16").Value = "TimeoutXL"
$constants.Range("B16").Value = 60000

# --- Restore the cursor / selection position on each sheet ---
$assets = $wb.Worksheets.Item("Assets")

$settings.Activate() | Out-Null
$settings.Range("C22").Select() | Out-Null

$assets.Activate() | Out-Null
$assets.Range("C25").Select() | Out-Null

$constants.Activate() | Out-Null
$constants.Range("B25").Select() | Out-Null
